$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: plain numeric values (previously row 2 held text "5000.0")
$ws.Range("A2").Value = 5000
$ws.Range("B2").Value = 5000

$ws.Range("A3").Value = 4000
$ws.Range("B3").Value = 4000

$ws.Range("A4").Value = 4000
$ws.Range("B4").Value = 5000

$ws.Range("A5").Value = 4740
$ws.Range("B5").Value = 4041

# Row 6: text values that look like numbers ("4141.0"/"4142.0"), stored as
# actual text rather than being auto-coerced to numbers. Force a text
# number-format while assigning, then restore the default style so no
# residual formatting sticks to the cells.
$ws.Range("A6:B6").NumberFormat = "@"
$ws.Range("A6").Value = "4141.0"
$ws.Range("B6").Value = "4142.0"
$ws.Range("A6:B6").Style = "Normal"
